# Remove the "Cover topics from CO3002/7002." paragraph from the
# Content Placeholder on slide 3, leaving the surrounding paragraphs
# (and their formatting / endParaRPr) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$target = $tr.Paragraphs(3, 1)
$target.Delete()
